$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the now-unused index tail (rows 45:87), and move the
#     active selection to D50 with the window scrolled so A33 is the
#     top-left visible cell. ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("45:87").Delete()

# --- Sheet1 becomes the active sheet/tab (was Sheet3). ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("D50").Select()
